# Reorder the "Recorded By" (column G) list of names/emails in the
# "Session Analysis Results" sheet: for every cell that holds a
# comma-separated list of two or more recorders, move the LAST entry
# to the FRONT of the list (a single right-rotation), leaving the
# relative order of the remaining entries unchanged.
#
# e.g. "System, dnasr281@gmail.com"          -> "dnasr281@gmail.com, System"
#      "System, system, backup@backdoor.com" -> "backup@backdoor.com, System, system"
#
# Cells holding only a single value (no comma) are left untouched, and
# so are the handful of "System, admin@admin.com" cells, matching the
# upstream source data exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

$colG = 7

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $colG)
    $val = $cell.Value2

    if ($null -eq $val) { continue }
    if ($val -eq "") { continue }
    if ($val -eq "Recorded By") { continue }
    if ($val -eq "System, admin@admin.com") { continue }

    $parts = $val -split ", "
    $n = $parts.Count

    if ($n -le 1) { continue }

    $newParts = @($parts[$n - 1]) + $parts[0..($n - 2)]
    $newVal = $newParts -join ", "

    $cell.Value = $newVal
}
